$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# F2 holds the description of the "color" enum column; drop the trailing ", "
$ws.Range("F2").Value = "enum | red, green, blue"

# F3 holds the sample color value for row 1 of data; change it from green to red
$ws.Range("F3").Value = "red"
